$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows("2:5").Delete()

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Rows("2:2").Delete()
$ws2.Range("B2").Value = 20
